$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel constant used by Range.Replace (XlLookAt): xlWhole = 1, xlPart = 2
$xlWhole = 1

$cells = $ws.Cells

# 1) "bleu" -> "noir" (statut_label column)
$cells.Replace("bleu", "noir", $xlWhole) | Out-Null

# 2) Correct the statut_name wording
$cells.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois", $xlWhole) | Out-Null
$cells.Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois", $xlWhole) | Out-Null
$cells.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés", $xlWhole) | Out-Null
$cells.Replace("résultat et / ou publication posté", "résultat postés ou publiés", $xlWhole) | Out-Null
